# edit.ps1 - apply the "render site after updating content" change:
#   1. Split the Title paragraph's single run into per-word/space runs.
#   2. Split the Question 6 curly-quoted phrase into three runs
#      (open-quote / phrase / close-quote).
#   3. Re-point the "Subtitle" paragraph style's basedOn from "Title" to "Normal".
#   4. Add an explicit theme color (text1 / tint A6) to the "Subtitle" style rPr.
#   5. Add an explicit color (345A8A) to the "AbstractTitle" style rPr.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the text that currently lives in a single run into several
# sibling runs (no formatting differences) by temporarily inserting paragraph
# breaks at the requested absolute character offsets and then deleting the
# inserted paragraph marks again. Word (and this COM emulation) keeps runs
# that originate from distinct paragraphs distinct even after the paragraphs
# are rejoined, whereas editing text in place inside one run/range collapses
# back into a single run.
# ---------------------------------------------------------------------------
function Split-RunAtOffsets {
    param(
        [int]$AnchorStart,
        [int[]]$Offsets
    )

    # Insert the breaks starting from the right-most offset so that earlier,
    # not-yet-processed offsets stay valid (inserting a paragraph mark shifts
    # everything after it by one character).
    for ($i = $Offsets.Length - 1; $i -ge 0; $i--) {
        $pos = $Offsets[$i]
        $breakRange = $d.Range($pos, $pos)
        $breakRange.InsertParagraphAfter() | Out-Null
    }

    # Now stitch the paragraphs back together by repeatedly deleting the
    # trailing paragraph mark of the paragraph that starts at $AnchorStart.
    for ($i = 0; $i -lt $Offsets.Length; $i++) {
        $para = $d.Range($AnchorStart, $AnchorStart).Paragraphs(1).Range
        $mark = $d.Range($para.End - 1, $para.End)
        $mark.Delete() | Out-Null
    }
}


# ---------------------------------------------------------------------------
# 1. Title: "Week 7 Day 2 – Bootstrap Activity" -> 13 separate runs.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titleStart = $titlePara.Start

$segments = @("Week", " ", "7", " ", "Day", " ", "2", " ", "–", " ", "Bootstrap", " ", "Activity")
$titleOffsets = @()
$cursor = $titleStart
for ($i = 0; $i -lt ($segments.Length - 1); $i++) {
    $cursor = $cursor + $segments[$i].Length
    $titleOffsets += $cursor
}

Split-RunAtOffsets $titleStart $titleOffsets

# ---------------------------------------------------------------------------
# 2. Question 6: “I am 95% confident,” -> “ / I am 95% confident, / ”
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("“I am 95% confident,”", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the Question 6 quoted phrase"
}
$quoteStart = $findRange.Start
$quoteEnd = $findRange.End
$quoteOffsets = @($quoteStart + 1, $quoteEnd - 1)

Split-RunAtOffsets $quoteStart $quoteOffsets

# ---------------------------------------------------------------------------
# 3 & 4. "Subtitle" paragraph style: basedOn Title -> Normal, plus an explicit
#         text1/A6 theme color (matching the already-present SubtitleChar).
# ---------------------------------------------------------------------------
$subtitleStyle = $d.Styles("Subtitle")
$normalStyle = $d.Styles("Normal")
$subtitleStyle.BaseStyle = $normalStyle

$subtitleCharStyle = $d.Styles("SubtitleChar")
$subtitleStyle.Font.Color = $subtitleCharStyle.Font.Color

# ---------------------------------------------------------------------------
# 5. "AbstractTitle" paragraph style: add an explicit color 345A8A.
# ---------------------------------------------------------------------------
$abstractTitleStyle = $d.Styles("AbstractTitle")
$abstractTitleStyle.Font.Color = 9067060   # RGB(0x34, 0x5A, 0x8A) -> 0x8A5A34

Write-Host "done"
